# Generate Report for Handoff
# Replaces the old localization-request id (197ddc51-a5cd-4e34-928c-42c744ac13b9)
# with the new one (078cca57-638b-4e66-8f52-ab84e25d647b) everywhere it is
# referenced, updates the regenerated handoff-package file names (new content
# hash) and bumps the handoff timestamps to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldId = "197ddc51-a5cd-4e34-928c-42c744ac13b9"
$newId = "078cca57-638b-4e66-8f52-ab84e25d647b"

$oldHashZh = "3160eb2bda84981c47afa894ad6f693e06e3bce7"
$newHashZh = "7735079850a38c5b42c70020eb6eaa5cbf1a0669"

$oldHashDe = "3160eb2bda84981c47afa894ad6f693e06e3bce7"
$newHashDe = "7735079850a38c5b42c70020eb6eaa5cbf1a0669"

$newZhTimestamp = "2016-03-03 13:02:51"
$newDeTimestamp = "2016-03-03 13:03:08"

$newMdName = "$newId.md"
$newZhXlfName = "$newId.$newHashZh.zh-cn.xlf"
$newDeXlfName = "$newId.$newHashDe.de-de.xlf"

function Update-Hyperlink($ws, $cellAddress, $newDisplay) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddress) {
            $hl.TextToDisplay = $newDisplay
        }
    }
}

# --- Sheet "Overview": A2 is the source .md file name/hyperlink ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
Update-Hyperlink $wsOverview '$A$2' $newMdName

# --- Sheet "zh-cn": A2 (.md), C2 (.xlf handoff file) + D2 (handoff datetime) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
Update-Hyperlink $wsZh '$A$2' $newMdName

$wsZh.Range("C2").Value = $newZhXlfName
Update-Hyperlink $wsZh '$C$2' $newZhXlfName

$wsZh.Range("D2").Value = $newZhTimestamp

# --- Sheet "de-de": A2 (.md), C2 (.xlf handoff file) + D2 (handoff datetime) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
Update-Hyperlink $wsDe '$A$2' $newMdName

$wsDe.Range("C2").Value = $newDeXlfName
Update-Hyperlink $wsDe '$C$2' $newDeXlfName

$wsDe.Range("D2").Value = $newDeTimestamp
